# Commit: "add cyclone_phailin and run topN, shifty on all events"
#
# The Cyclone in India ("cyclone_phailin") event row (row 17) gets its
# DataFile filename filled in, and its Keywords text corrected from
# "gopalpur" to "odisha". The trailing, always-empty row (row 22) is
# removed from the event log. Finally the active selection is left on
# A20, matching where the author's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Keywords (column F) updated first, then DataFile (column B) - this
# mirrors the order new shared-string entries were appended in the
# original edit.
$ws.Range("F17").Value = "cyclone, phailin, odisha"
$ws.Range("B17").Value = "2013_10_11_20"

# Drop the blank trailing row at the bottom of the table.
$ws.Rows.Item(22).Delete()

# Leave the selection where the author left it.
$ws.Range("A20").Select() | Out-Null
